$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 5) to the table, e.g. series with date
# "11-08-2021". A plain string assignment of a date-looking value would
# be auto-converted by Excel into a date serial number (with an
# associated date/number-format style), which does not match the
# source data (a plain text shared string, like the other date-looking
# labels in column A). To force text storage without leaving a stray
# cell style behind, compute the text via a formula and then convert
# the formula result to a static value with copy/paste-special.
$ws.Range("A5").Formula = '=TEXT("11-08-2021","@")'
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4163) | Out-Null   # xlPasteValues

$ws.Range("B5").Value = 3500
$ws.Range("C5").Value = 9520
$ws.Range("D5").Value = 3500
$ws.Range("E5").Value = 50
$ws.Range("F5").Value = 3450
$ws.Range("G5").Value = 2.02
